$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (VFTS 243 / Cyg X-1 row) citation info gained a second reference
# (a 2012 ApJ paper on 1107.5585) alongside the existing 2021 Science paper.
# Set the link (K) before the label (J) so that, in the shared-string table,
# the URL string ends up with the lower index and the citation text with the
# higher index - matching the ordering produced by the source edit.
$ws.Range("K3").Value = "https://ui.adsabs.harvard.edu/abs/2021Sci...371.1046M/abstract, https://arxiv.org/abs/1107.5585"
$ws.Range("J3").Value = "2021Sci...371.1046M , 2012ApJ...747..111W"

# Update the sheet's active selection to G3 (was J13)
$ws.Range("G3").Select()
